$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SampleInfo")
$ws2 = $wb.Worksheets.Item("SampleRenames")

# --- SampleInfo sheet: append "_IGO" to CMO_SAMPLE_ID (col A), and set
#     INVESTIGATOR_SAMPLE_ID (col C) to the clean investigator-style id ---
$ws1.Range("A2").Value = "test_sample_2_T_IGO"
$ws1.Range("C2").Value = "test_investigator_sample_2_T"

$ws1.Range("A3").Value = "test_sample_1_N_IGO"
$ws1.Range("C3").Value = "test_investigator_sample_1_N"

$ws1.Range("A4").Value = "test_sample_4_T_IGO"
$ws1.Range("C4").Value = "test_investigator_sample_4_T"

$ws1.Range("A5").Value = "test_sample_3_N_IGO"
$ws1.Range("C5").Value = "test_investigator_sample_3_N"

$ws1.Range("A6").Value = "test_sample_6_T_IGO"
$ws1.Range("C6").Value = "test_investigator_sample_6_T"

$ws1.Range("A7").Value = "test_sample_5_N_IGO"
$ws1.Range("C7").Value = "test_investigator_sample_5_N"

# --- SampleRenames sheet: OldName (col A) gets "_IGO" suffix, NewName
#     (col B) becomes the clean sample id (drop test-only suffixes) ---
$ws2.Range("A2").Value = "test_sample_2_T_IGO"
$ws2.Range("B2").Value = "test_sample_2_T"

$ws2.Range("A3").Value = "test_sample_1_N_IGO"
$ws2.Range("B3").Value = "test_sample_1_N"

$ws2.Range("A4").Value = "test_sample_4_T_IGO"

$ws2.Range("A5").Value = "test_sample_3_N_IGO"

$ws2.Range("A6").Value = "test_sample_6_T_IGO"
$ws2.Range("B6").Value = "test_sample_6_T"

$ws2.Range("A7").Value = "test_sample_5_N_IGO"

# Column A on SampleRenames drops its explicit style (back to Normal/default)
# while column B keeps its existing style.
$ws2.Range("A2:A7").Style = "Normal"

# --- Selections, matching the saved view state in each sheet ---
$ws1.Range("I11").Select()
$ws2.Range("C25").Select()
$ws2.Activate()
